# Update gh-pages to output generated at 456a3b4
# Applies refreshed "想去人数" (interested-people count) figures across the
# four worksheets, matching the source data regeneration.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1260
$ws.Range("F5").Value  = 5555
$ws.Range("F6").Value  = 1785
$ws.Range("F7").Value  = 6351
$ws.Range("F8").Value  = 140
$ws.Range("F9").Value  = 1914
$ws.Range("F10").Value = 513
$ws.Range("F11").Value = 6
$ws.Range("F13").Value = 30
$ws.Range("F17").Value = 7900
$ws.Range("F18").Value = 7900
$ws.Range("F21").Value = 182
$ws.Range("F22").Value = 108
$ws.Range("F23").Value = 1744
$ws.Range("F30").Value = 1743
$ws.Range("F31").Value = 801
$ws.Range("F32").Value = 372
$ws.Range("F35").Value = 79
$ws.Range("F37").Value = 3916

# --- 演出 sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 361
$ws.Range("F5").Value  = 204
$ws.Range("F12").Value = 9
$ws.Range("F14").Value = 26

# --- 本地生活 sheet ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2274
$ws.Range("F4").Value = 683
$ws.Range("F5").Value = 266

# --- 全部类型 sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 2274
$ws.Range("F4").Value  = 683
$ws.Range("F5").Value  = 1260
$ws.Range("F9").Value  = 361
$ws.Range("F10").Value = 5555
$ws.Range("F11").Value = 266
$ws.Range("F12").Value = 1784
$ws.Range("F13").Value = 6351
$ws.Range("F14").Value = 140
$ws.Range("F15").Value = 1914
$ws.Range("F17").Value = 513
$ws.Range("F19").Value = 30
$ws.Range("F23").Value = 7900
$ws.Range("F24").Value = 7900
$ws.Range("F27").Value = 182
$ws.Range("F28").Value = 108
$ws.Range("F29").Value = 1744
$ws.Range("F35").Value = 1743
$ws.Range("F36").Value = 801
$ws.Range("F38").Value = 372
$ws.Range("F40").Value = 26
$ws.Range("F45").Value = 3916
